$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("B2") 'Bitcoin'
Set-TextValue $ws.Range("C2") 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextValue $ws.Range("D2") '95.236.24'
Set-TextValue $ws.Range("E2") '  -1.84%  '
Set-TextValue $ws.Range("B3") 'Ethereum'
Set-TextValue $ws.Range("C3") 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextValue $ws.Range("D3") '3.603.04'
Set-TextValue $ws.Range("E3") '  -2.36%  '
Set-TextValue $ws.Range("B4") 'XRP'
Set-TextValue $ws.Range("C4") 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue $ws.Range("D4") '2.37'
Set-TextValue $ws.Range("E4") '  +24.29%  '
Set-TextValue $ws.Range("B5") 'TetherUSD'
Set-TextValue $ws.Range("C5") 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextValue $ws.Range("D5") '0.998'
Set-TextValue $ws.Range("E5") '  -0.26%  '
Set-TextValue $ws.Range("B6") 'Solana'
Set-TextValue $ws.Range("C6") 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Range("D6") '225.02'
Set-TextValue $ws.Range("E6") '  -5.14%  '
Set-TextValue $ws.Range("B7") 'BNB'
Set-TextValue $ws.Range("C7") 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue $ws.Range("D7") '636.11'
Set-TextValue $ws.Range("E7") '  -3.21%  '
Set-TextValue $ws.Range("B8") 'Dogecoin'
Set-TextValue $ws.Range("C8") 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range("D8") '0.413'
Set-TextValue $ws.Range("E8") '  -2.65%  '
Set-TextValue $ws.Range("B9") 'Cardano'
Set-TextValue $ws.Range("C9") 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Range("D9") '1.09'
Set-TextValue $ws.Range("E9") '  +2.24%  '
Set-TextValue $ws.Range("B10") 'USDC'
Set-TextValue $ws.Range("C10") 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue $ws.Range("D10") '0.999'
Set-TextValue $ws.Range("E10") '  -0.05%  '
Set-TextValue $ws.Range("B11") 'LidoStakedEther'
Set-TextValue $ws.Range("C11") 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
Set-TextValue $ws.Range("D11") '3.601.34'
Set-TextValue $ws.Range("E11") '  -2.36%  '
Set-TextValue $ws.Range("B12") 'Avalanche'
Set-TextValue $ws.Range("C12") 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range("D12") '46.87'
Set-TextValue $ws.Range("E12") '  +6.51%  '
Set-TextValue $ws.Range("B13") 'TRON'
Set-TextValue $ws.Range("C13") 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D13") '0.207'
Set-TextValue $ws.Range("E13") '  -1.16%  '
Set-TextValue $ws.Range("B14") 'ShibaInu'
Set-TextValue $ws.Range("C14") 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D14") '0.0000289'
Set-TextValue $ws.Range("E14") '  -7.93%  '
Set-TextValue $ws.Range("B15") 'Toncoin'
Set-TextValue $ws.Range("C15") 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D15") '6.47'
Set-TextValue $ws.Range("E15") '  -4.22%  '
Set-TextValue $ws.Range("B16") 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range("C16") 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range("D16") '4.274.08'
Set-TextValue $ws.Range("E16") '  -2.35%  '
Set-TextValue $ws.Range("B17") 'WrappedBTC'
Set-TextValue $ws.Range("C17") 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range("D17") '94.925.02'
Set-TextValue $ws.Range("E17") '  -1.87%  '
Set-TextValue $ws.Range("B18") 'Polkadot'
Set-TextValue $ws.Range("C18") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D18") '8.75'
Set-TextValue $ws.Range("E18") '  -4.16%  '
Set-TextValue $ws.Range("B19") 'WrappedEther'
Set-TextValue $ws.Range("C19") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D19") '3.601.62'
Set-TextValue $ws.Range("E19") '  -2.58%  '
Set-TextValue $ws.Range("B20") 'Uniswap'
Set-TextValue $ws.Range("C20") 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D20") '13.45'
Set-TextValue $ws.Range("E20") '  +3.57%  '
Set-TextValue $ws.Range("B21") 'Chainlink'
Set-TextValue $ws.Range("C21") 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D21") '19.63'
Set-TextValue $ws.Range("E21") '  +5.23%  '
Set-TextValue $ws.Range("B22") 'Stellar'
Set-TextValue $ws.Range("C22") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D22") '0.516'
Set-TextValue $ws.Range("E22") '  +1.84%  '
Set-TextValue $ws.Range("B23") 'BitcoinCash'
Set-TextValue $ws.Range("C23") 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D23") '501.01'
Set-TextValue $ws.Range("E23") '  -3.67%  '
Set-TextValue $ws.Range("B24") 'SuiNetwork'
Set-TextValue $ws.Range("C24") 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range("D24") '3.23'
Set-TextValue $ws.Range("E24") '  -6.47%  '
Set-TextValue $ws.Range("B25") 'Hedera'
Set-TextValue $ws.Range("C25") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D25") '0.239'
Set-TextValue $ws.Range("E25") '  +21.34%  '
Set-TextValue $ws.Range("B26") 'Litecoin'
Set-TextValue $ws.Range("C26") 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D26") '120.22'
Set-TextValue $ws.Range("E26") '  +18.67%  '
Set-TextValue $ws.Range("B27") 'PEPE'
Set-TextValue $ws.Range("C27") 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D27") '0.0000202'
Set-TextValue $ws.Range("E27") '  -4.96%  '
Set-TextValue $ws.Range("B28") 'NEARProtocol'
Set-TextValue $ws.Range("C28") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D28") '6.73'
Set-TextValue $ws.Range("E28") '  -3.01%  '
Set-TextValue $ws.Range("B29") 'WrappedeETH'
Set-TextValue $ws.Range("C29") 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue $ws.Range("D29") '3.793.30'
Set-TextValue $ws.Range("E29") '  -2.47%  '
Set-TextValue $ws.Range("B30") 'Aptos'
Set-TextValue $ws.Range("C30") 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D30") '12.60'
Set-TextValue $ws.Range("E30") '  -6.33%  '
Set-TextValue $ws.Range("B31") 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D31") '12.91'
Set-TextValue $ws.Range("E31") '  +3.31%  '
Set-TextValue $ws.Range("B32") 'PancakeSwap'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D32") '2.92'
Set-TextValue $ws.Range("E32") '  -3.50%  '
Set-TextValue $ws.Range("B33") 'Dai'
Set-TextValue $ws.Range("C33") 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D33") '1.00'
Set-TextValue $ws.Range("E33") '  +0.03%  '
Set-TextValue $ws.Range("B34") 'Binance-PegBSC-USD'
Set-TextValue $ws.Range("C34") 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D34") '1.00'
Set-TextValue $ws.Range("E34") '  -0.12%  '
Set-TextValue $ws.Range("B35") 'Cronos'
Set-TextValue $ws.Range("C35") 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D35") '0.178'
Set-TextValue $ws.Range("E35") '  -5.44%  '
Set-TextValue $ws.Range("B36") 'EthereumClassic'
Set-TextValue $ws.Range("C36") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D36") '31.92'
Set-TextValue $ws.Range("E36") '  -0.58%  '
Set-TextValue $ws.Range("B37") 'Fetch.AI'
Set-TextValue $ws.Range("C37") 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D37") '1.76'
Set-TextValue $ws.Range("E37") '  -5.99%  '
Set-TextValue $ws.Range("B38") 'PolygonEcosystemToken'
Set-TextValue $ws.Range("C38") 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range("D38") '0.588'
Set-TextValue $ws.Range("E38") '  -0.36%  '
Set-TextValue $ws.Range("B39") 'USDe'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D39") '1.00'
Set-TextValue $ws.Range("E39") '  +0.01%  '
Set-TextValue $ws.Range("B40") 'Bittensor'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D40") '591.99'
Set-TextValue $ws.Range("E40") '  -9.32%  '
Set-TextValue $ws.Range("B41") 'RenderToken'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range("D41") '8.30'
Set-TextValue $ws.Range("E41") '  -5.61%  '
Set-TextValue $ws.Range("B42") 'EnergySwap'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D42") '41.95'
Set-TextValue $ws.Range("E42") '  +3.44%  '
Set-TextValue $ws.Range("B43") 'Filecoin'
Set-TextValue $ws.Range("C43") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D43") '6.85'
Set-TextValue $ws.Range("E43") '  +0.95%  '
Set-TextValue $ws.Range("B44") 'Kaspa'
Set-TextValue $ws.Range("C44") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D44") '0.158'
Set-TextValue $ws.Range("E44") '  -2.12%  '
Set-TextValue $ws.Range("B45") 'Algorand'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D45") '0.479'
Set-TextValue $ws.Range("E45") '  -0.87%  '
Set-TextValue $ws.Range("B46") 'ImmutableX'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D46") '1.91'
Set-TextValue $ws.Range("E46") '  -6.00%  '
Set-TextValue $ws.Range("B47") 'VeChain'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D47") '0.0472'
Set-TextValue $ws.Range("E47") '  +1.63%  '
Set-TextValue $ws.Range("B48") 'ARBITRUM'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D48") '0.919'
Set-TextValue $ws.Range("E48") '  -4.49%  '
Set-TextValue $ws.Range("B49") 'Aave'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D49") '225.78'
Set-TextValue $ws.Range("E49") '  +10.50%  '
Set-TextValue $ws.Range("B50") 'WhiteBITCoin'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range("D50") '23.47'
Set-TextValue $ws.Range("E50") '  -0.67%  '
Set-TextValue $ws.Range("B51") 'MantraDAO'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue $ws.Range("D51") '3.67'
Set-TextValue $ws.Range("E51") '  +5.04%  '
